$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "65.449.96"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +2.99%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.659.50"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +2.24%  "

$ws.Range("E4").Value = "  -0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "606.62"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.05%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "157.24"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +3.68%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +1.18%  "

$ws.Range("E9").Value = "  +7.41%  "

$ws.Range("E10").Value = "  +4.37%  "

$ws.Range("E11").Value = "  +3.65%  "

$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("E13").Value = "  +5.86%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "3.135.52"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +2.13%  "

$ws.Range("E15").Value = "  +15.91%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "65.305.93"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.95%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.874.49"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +10.39%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "12.67"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "

$ws.Range("E19").Value = "  +2.45%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "355.75"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +2.41%  "

$ws.Range("E21").Value = "  +5.82%  "

$ws.Range("E22").Value = "  +0.17%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "68.42"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("E25").Value = "  +2.65%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "8.33"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +3.67%  "

$ws.Range("E28").Value = "  +2.26%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "548.38"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -3.83%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0958"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +12.26%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("E32").Value = "  +3.20%  "

$ws.Range("E33").Value = "  +3.41%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.80"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +10.97%  "

$ws.Range("E35").Value = "  +4.94%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.430"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +3.83%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.07"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +6.28%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "165.64"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.00%  "

$ws.Range("E39").Value = "  +3.23%  "

$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("E41").Value = "  -0.01%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "168.34"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.23%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "42.37"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +6.87%  "

$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("E45").Value = "  +4.39%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "23.41"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +4.48%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.26"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +9.70%  "

$ws.Range("E48").Value = "  +4.13%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0254"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("E50").Value = "  +1.94%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "19.77"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.94%  "

